$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-10 from 2023-10-22 (45221)
# to 2023-10-25 (45224), preserving existing date formatting/style.
for ($row = 2; $row -le 10; $row++) {
    $ws.Range("C$row").Value = 45224
}
